# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below would be auto-coerced from text to numbers by Excel because their
# new values look like plain decimals. Mark them as Text format first so the
# values are stored as strings, matching the source data (inline strings).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Now write the updated values / text cells.
$ws.Range('D2').Value = '44.093.26'
$ws.Range('E2').Value = '  +5.68%  '
$ws.Range('D3').Value = '2.278.96'
$ws.Range('E3').Value = '  +3.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '232.94'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '0.643'
$ws.Range('E6').Value = '  +4.08%  '
$ws.Range('E7').Value = '  +9.16%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '0.437'
$ws.Range('E9').Value = '  +8.85%  '
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  +16.08%  '
$ws.Range('D11').Value = '57.55'
$ws.Range('E11').Value = '  +0.46%  '
$ws.Range('D12').Value = '26.23'
$ws.Range('E12').Value = '  +18.54%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '2.616.25'
$ws.Range('D15').Value = '15.98'
$ws.Range('E15').Value = '  +3.86%  '
$ws.Range('D16').Value = '6.02'
$ws.Range('E16').Value = '  +7.75%  '
$ws.Range('D17').Value = '0.838'
$ws.Range('E17').Value = '  +5.59%  '
$ws.Range('D18').Value = '2.276.36'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('D19').Value = '43.951.89'
$ws.Range('E19').Value = '  +5.49%  '
$ws.Range('D20').Value = '0.0₃0991'
$ws.Range('E20').Value = '  +10.21%  '
$ws.Range('D21').Value = '74.05'
$ws.Range('E21').Value = '  +2.89%  '
$ws.Range('D22').Value = '6.14'
$ws.Range('D23').Value = '260.33'
$ws.Range('E23').Value = '  +7.59%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '2.49'
$ws.Range('E25').Value = '  +6.11%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '2.32'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '10.22'
$ws.Range('E27').Value = '  +5.95%  '
$ws.Range('D28').Value = '171.86'
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('D29').Value = '21.04'
$ws.Range('E29').Value = '  +6.58%  '
$ws.Range('D30').Value = '0.139'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('E32').Value = '  +7.29%  '
$ws.Range('E33').Value = '  +2.34%  '
$ws.Range('D34').Value = '0.0692'
$ws.Range('E34').Value = '  +7.11%  '
$ws.Range('D35').Value = '5.04'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').Value = '4.75'
$ws.Range('E36').Value = '  +2.82%  '
$ws.Range('D37').Value = '3.88'
$ws.Range('E37').Value = '  +9.53%  '
$ws.Range('D38').Value = '6.84'
$ws.Range('E38').Value = '  +8.50%  '
$ws.Range('D39').Value = '2.37'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').Value = '0.0249'
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '8.39'
$ws.Range('E42').Value = '  -2.31%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '17.61'
$ws.Range('E43').Value = '  +8.58%  '
$ws.Range('D44').Value = '0.0976'
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').Value = '4.49'
$ws.Range('E45').Value = '  +1.91%  '
$ws.Range('D46').Value = '98.14'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '1.20'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  +7.04%  '
$ws.Range('D49').Value = '1.472.81'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').Value = '9.91'
$ws.Range('E50').Value = '  +17.37%  '
$ws.Range('D51').Value = '0.000203'
$ws.Range('E51').Value = '  -14.84%  '
